{"js": "// The document is a single-column table where each row holds one\n// benchmark statistic. This edit updates the text of a fixed set of\n// rows (by position) to their new values, matching the target diff:\n//  - rows 0-2 become \"0M\"\n//  - row 3 becomes \"33\" (the following three legacy rows are folded away)\n//  - rows 4-11 get refreshed numeric values (row 7 is unchanged)\n//  - the three \"summary\" rows (formerly multi-column lines separated by\n//    tabs) collapse down to their single headline percentage value.\nconst changes = [\n  { idx: 0, value: \"0M\" },\n  { idx: 1, value: \"0M\" },\n  { idx: 2, value: \"0M\" },\n  { idx: 3, value: \"33\" },\n  { idx: 4, value: \"0.00004\" },\n  { idx: 5, value: \"0.00055\" },\n  { idx: 6, value: \"0.00020\" },\n  { idx: 8, value: \"0.00041\" },\n  { idx: 9, value: \"0.00045\" },\n  { idx: 10, value: \"0.00051\" },\n  { idx: 11, value: \"0.00808\" },\n  { idx: 43, value: \"99.99\" },\n  { idx: 44, value: \"0.01\" },\n  { idx: 45, value: \"65\" },\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rows = table.rows.items;\nconst maxIdx = Math.max(...changes.map((c) => c.idx));\nfor (let i = 0; i <= maxIdx && i < rows.length; i++) {\n  rows[i].cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const { idx, value } of changes) {\n  if (idx < rows.length) {\n    rows[idx].cells.items[0].value = value;\n  }\n}\nawait context.sync();\n", "ps1": "# The document is a single-column table where each row holds one\n# benchmark statistic. This edit updates the text of a fixed set of\n# rows (by position, 1-based as COM indexes) to their new values,\n# matching the target diff:\n#  - rows 1-3 become \"0M\"\n#  - row 4 becomes \"33\" (the following three legacy rows are folded away)\n#  - rows 5-12 get refreshed numeric values (row 8 is unchanged)\n#  - the three \"summary\" rows (formerly multi-column lines separated by\n#    tabs) collapse down to their single headline percentage value.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$changes = @(\n    @{ Row = 1; Value = \"0M\" },\n    @{ Row = 2; Value = \"0M\" },\n    @{ Row = 3; Value = \"0M\" },\n    @{ Row = 4; Value = \"33\" },\n    @{ Row = 5; Value = \"0.00004\" },\n    @{ Row = 6; Value = \"0.00055\" },\n    @{ Row = 7; Value = \"0.00020\" },\n    @{ Row = 9; Value = \"0.00041\" },\n    @{ Row = 10; Value = \"0.00045\" },\n    @{ Row = 11; Value = \"0.00051\" },\n    @{ Row = 12; Value = \"0.00808\" },\n    @{ Row = 44; Value = \"99.99\" },\n    @{ Row = 45; Value = \"0.01\" },\n    @{ Row = 46; Value = \"65\" }\n)\n\nforeach ($change in $changes) {\n    $t.Cell($change.Row, 1).Range.Text = $change.Value\n}\n"}
